$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.029.06'
$ws.Range("E2").Value = '  +0.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.679.33'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.15'
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("E6").Value = '  -2.86%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.61'
$ws.Range("E8").Value = '  +7.18%  '

$ws.Range("E9").Value = '  +1.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0622'
$ws.Range("E10").Value = '  +0.62%  '

$ws.Range("E11").Value = '  -0.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.916.40'
$ws.Range("E12").Value = '  +0.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.651.39'
$ws.Range("E13").Value = '  -0.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.532'
$ws.Range("E15").Value = '  +1.66%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.44'
$ws.Range("E16").Value = '  +0.76%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.039.96'
$ws.Range("E17").Value = '  +0.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.20'
$ws.Range("E18").Value = '  +5.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.54'
$ws.Range("E19").Value = '  +1.60%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0739'
$ws.Range("E20").Value = '  +1.01%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.46'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("E24").Value = '  -4.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.62'
$ws.Range("E25").Value = '  +0.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.24'
$ws.Range("E26").Value = '  +1.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.49'
$ws.Range("E27").Value = '  +3.84%  '

$ws.Range("E28").Value = '  -2.52%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  +0.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.517.54'
$ws.Range("E33").Value = '  +3.42%  '

$ws.Range("E34").Value = '  +1.35%  '

$ws.Range("E35").Value = '  +4.76%  '

$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.591'
$ws.Range("E37").Value = '  +3.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.918'
$ws.Range("E38").Value = '  +2.36%  '

$ws.Range("E39").Value = '  +3.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  +5.53%  '

$ws.Range("E41").Value = '  -2.16%  '

$ws.Range("E42").Value = '  +0.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '67.91'
$ws.Range("E43").Value = '  +3.25%  '

$ws.Range("E44").Value = '  -0.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.821.98'
$ws.Range("E45").Value = '  +0.38%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.19'
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.53'
$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("E49").Value = '  +3.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.95'
$ws.Range("E50").Value = '  +4.85%  '

$ws.Range("E51").Value = '  -0.56%  '
